$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -8.250999999999999
$ws.Range("B4").Value = 8.188999999999998
$ws.Range("D4").Value = -7.830000000000001
$ws.Range("D5").Value = -8.248000000000001
$ws.Range("A8").Value = -21.699
$ws.Range("A10").Value = -21.808
$ws.Range("A12").Value = -21.303
$ws.Range("B12").Value = 6.303999999999999
$ws.Range("C12").Value = -11.263
$ws.Range("C13").Value = -12.813
$ws.Range("B15").Value = 5.394
$ws.Range("B17").Value = 4.836999999999999
$ws.Range("A18").Value = -21.734
$ws.Range("D20").Value = -7.879
$ws.Range("C21").Value = -12.808
$ws.Range("D23").Value = -8.280999999999999
$ws.Range("C25").Value = -12.151
$ws.Range("B26").Value = 5.525999999999999
$ws.Range("D26").Value = -7.669000000000001
$ws.Range("B27").Value = 5.659000000000001
$ws.Range("B28").Value = 5.119
$ws.Range("C32").Value = -12.439
$ws.Range("D34").Value = -7.74
$ws.Range("C36").Value = -12.776
$ws.Range("A37").Value = -21.2
$ws.Range("B37").Value = 6.423999999999999
$ws.Range("C38").Value = -12.356
$ws.Range("D39").Value = -7.476000000000001
$ws.Range("D40").Value = -7.922
$ws.Range("C41").Value = -12.775
$ws.Range("D41").Value = -7.741
$ws.Range("B47").Value = 5.435
$ws.Range("D47").Value = -8.057
$ws.Range("C52").Value = -11.775
$ws.Range("A55").Value = -22.109
$ws.Range("C59").Value = -12.7
$ws.Range("D60").Value = -8.321
$ws.Range("B65").Value = 5.522
$ws.Range("C67").Value = -10.944
$ws.Range("A68").Value = -21.502
$ws.Range("D72").Value = -7.525
$ws.Range("B73").Value = 6.792
$ws.Range("A77").Value = -20.79
$ws.Range("A78").Value = -20.425
$ws.Range("A81").Value = -21.747
$ws.Range("A82").Value = -21.822
$ws.Range("D83").Value = -7.944
$ws.Range("B84").Value = 4.991000000000001
$ws.Range("C84").Value = -12.819
$ws.Range("B85").Value = 5.008
$ws.Range("C88").Value = -13.222
$ws.Range("C89").Value = -13.791
$ws.Range("B93").Value = 5.587000000000001
$ws.Range("B95").Value = 6.332
$ws.Range("C95").Value = -12.248
$ws.Range("B98").Value = 7.208
$ws.Range("B99").Value = 5.285000000000001
$ws.Range("B101").Value = 6.043000000000001
$ws.Range("C105").Value = -12.753
